$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B13").Value = "Entspanntes Jahresgespräch"
$ws.Range("E13").Value = "https://cdntest.apozin.de/mimg/archive/large/adobestock_284485921_entspannung.jpg"
$ws.Range("D13").Value = "Bei mir sein"
$ws.Range("C13").Value = "Ich hatte mein Jahresgespräch mit meinem Vorgesetzen. Ich war maximal entspannt. In den Jahren zuvor hatte ich immer das Gefühl, ich müsse die Veränderung anstossen bzw. kritisieren. Ich erhoffte mir dann, dass meine Vorstellung durchkommt und dadurch kam ich ins emotionale Ungleichgewicht. Diesmal konnte ich meine Energie bei mir halten und war auch sachlich und entspannt bei der Lohndiskussion. "

$ws.Range("D15").Select()
